# rfe: sort results by package name
#
# The report's results table (Sheet1!A5:J10 - header in row 5, the
# "Package Name" column is A) needs its data rows sorted alphabetically
# by Package Name, exactly like Excel's Data > Sort would do: the whole
# row (all columns, plus per-cell formatting such as the orange
# "Has invalid ..." highlight) moves together with its Package Name key.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlAscending = 1
$xlYes = 1

# Table range including the header row, and the sort key (Package Name column).
$dataRange = $ws.Range("A5:J10")
$keyRange = $ws.Range("A5:A10")

$dataRange.Sort($keyRange, $xlAscending, $null, $null, $null, $null, $null, $xlYes)
